$d = $word.ActiveDocument

$replacements = @(
    @("85÷7=", "18÷7="),
    @("76÷2=", "54÷6="),
    @("96÷3=", "76÷3="),
    @("82÷7=", "45÷3="),
    @("61÷2=", "78÷6="),
    @("62÷9=", "65÷8="),
    @("39÷9=", "11÷6="),
    @("99÷3=", "96÷2="),
    @("41÷4=", "69÷2="),
    @("28÷4=", "24÷7="),
    @("59÷2=", "76÷5="),
    @("55÷7=", "20÷2="),
    @("88÷6=", "24÷4="),
    @("60÷7=", "40÷4="),
    @("15÷2=", "14÷6="),
    @("78÷7=", "59÷6="),
    @("73÷3=", "17÷6="),
    @("74÷8=", "38÷8="),
    @("24÷5=", "87÷9="),
    @("91÷2=", "14÷3="),
    @("27÷4=", "27÷5="),
    @("92÷9=", "97÷4="),
    @("98÷5=", "71÷7="),
    @("85÷2=", "20÷9="),
    @("76÷6=", "71÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
